$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix groupsize for the PHQ/HADS (Fischer et al.) row: E13 was left blank,
# should be "194". Enter with a leading apostrophe so Excel stores it as
# text (matching the other groupsize cells in column E, which are all
# text values) instead of auto-converting it to a number.
$cell = $ws.Range("E13")
$cell.Value = "'194"
